# Refactored portfolio to use instrument
# Insert a new "Instrument" column before the existing "Sector" column (J),
# shifting Sector/Category/Sub Category/Startup/Investment Domicile one
# column to the right, and populate the new column with "Stock".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J currently holds "Sector" - insert a new blank column before it,
# shifting Sector (and everything to its right) one column over.
$ws.Range("J1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("J1").Value = "Instrument"

# Data rows get the new "Stock" value.
$ws.Range("J2").Value = "Stock"
$ws.Range("J3").Value = "Stock"

# Restore the selection to match the post-edit state.
$ws.Range("J4").Select()
